$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("AR3").Value = 5
$ws.Range("AS3").Value = 1.17
# Row 5
$ws.Range("G5").Value = 2.3
$ws.Range("H5").Value = 2.9
$ws.Range("I5").Value = 3.6
$ws.Range("AF5").Value = 5.5
$ws.Range("AJ5").Value = 8
$ws.Range("AK5").Value = 17
$ws.Range("AL5").Value = 15
$ws.Range("AN5").Value = 41
# Row 6
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 7
$ws.Range("Q6").Value = 2.6
$ws.Range("R6").Value = 1.48
$ws.Range("S6").Value = 5.5
$ws.Range("T6").Value = 1.14
$ws.Range("AP6").Value = 1.98
$ws.Range("AQ6").Value = 1.88
$ws.Range("AR6").Value = 4.2
$ws.Range("AS6").Value = 1.23
# Row 8
$ws.Range("G8").Value = 1.91
$ws.Range("I8").Value = 3.2
$ws.Range("O8").Value = 1.18
$ws.Range("P8").Value = 4.5
$ws.Range("Q8").Value = 1.62
$ws.Range("R8").Value = 2.25
$ws.Range("S8").Value = 2.5
$ws.Range("T8").Value = 1.5
# Row 9
$ws.Range("G9").Value = 2.5
$ws.Range("I9").Value = 2.45
$ws.Range("J9").Value = 3.25
$ws.Range("O9").Value = 1.25
$ws.Range("P9").Value = 3.75
$ws.Range("Q9").Value = 1.88
$ws.Range("R9").Value = 1.98
$ws.Range("W9").Value = 1.67
$ws.Range("X9").Value = 2.1
$ws.Range("Y9").Value = 9.5
$ws.Range("AG9").Value = 13
$ws.Range("AJ9").Value = 9.5
$ws.Range("AK9").Value = 13
$ws.Range("AM9").Value = 26
# Row 10
$ws.Range("G10").Value = 3.25
$ws.Range("H10").Value = 3.6
$ws.Range("K10").Value = 2.3
$ws.Range("M10").Value = 1.03
$ws.Range("N10").Value = 15
$ws.Range("O10").Value = 1.2
$ws.Range("P10").Value = 4.33
$ws.Range("Q10").Value = 1.67
$ws.Range("R10").Value = 2.15
$ws.Range("S10").Value = 2.63
$ws.Range("T10").Value = 1.44
$ws.Range("W10").Value = 1.57
$ws.Range("X10").Value = 2.25
$ws.Range("Y10").Value = 13
$ws.Range("AB10").Value = 34
$ws.Range("AC10").Value = 23
$ws.Range("AE10").Value = 15
$ws.Range("AF10").Value = 7.5
$ws.Range("AJ10").Value = 9.5
$ws.Range("AO10").Value = 21
# Row 11
$ws.Range("G11").Value = 1.7
$ws.Range("H11").Value = 3.5
$ws.Range("I11").Value = 5.25
$ws.Range("N11").Value = 7.5
$ws.Range("Q11").Value = 2.25
$ws.Range("R11").Value = 1.62
$ws.Range("S11").Value = 4.33
$ws.Range("T11").Value = 1.2
$ws.Range("Z11").Value = 7
$ws.Range("AC11").Value = 17
$ws.Range("AE11").Value = 7.5
$ws.Range("AL11").Value = 17
$ws.Range("AN11").Value = 41
$ws.Range("AP11").Value = 1.71
$ws.Range("AQ11").Value = 2.11
$ws.Range("AR11").Value = 3.45
$ws.Range("AS11").Value = 1.3
# Row 13
$ws.Range("G13").Value = 2.4
$ws.Range("H13").Value = 2.7
$ws.Range("I13").Value = 3.4
$ws.Range("J13").Value = 3.15
$ws.Range("L13").Value = 4.05
$ws.Range("M13").Value = 1.13
$ws.Range("N13").Value = 5.1
$ws.Range("O13").Value = 1.53
$ws.Range("P13").Value = 2.35
$ws.Range("Q13").Value = 2.55
$ws.Range("R13").Value = 1.45
$ws.Range("S13").Value = 4.55
$ws.Range("T13").Value = 1.16
$ws.Range("W13").Value = 2.02
$ws.Range("AE13").Value = 5.1
$ws.Range("AF13").Value = 5.4
$ws.Range("AG13").Value = 16
$ws.Range("AH13").Value = 100
$ws.Range("AJ13").Value = 7.8
$ws.Range("AK13").Value = 17
$ws.Range("AM13").Value = 50
$ws.Range("AN13").Value = 37
$ws.Range("AO13").Value = 50
# Row 14
$ws.Range("G14").Value = 2.25
$ws.Range("I14").Value = 3.6
$ws.Range("J14").Value = 3.1
$ws.Range("L14").Value = 4.33
$ws.Range("Z14").Value = 9.5
$ws.Range("AB14").Value = 21
$ws.Range("AE14").Value = 6
$ws.Range("AJ14").Value = 8
$ws.Range("AO14").Value = 51
# Row 15
$ws.Range("G15").Value = 1.38
$ws.Range("I15").Value = 7
$ws.Range("AI15").Value = 351
$ws.Range("AK15").Value = 41
$ws.Range("AL15").Value = 21
# Row 16
$ws.Range("Q16").Value = 1.93
$ws.Range("R16").Value = 1.93
# Row 17
$ws.Range("M17").Value = 1.04
$ws.Range("N17").Value = 13
$ws.Range("Q17").Value = 1.83
$ws.Range("R17").Value = 2.03
$ws.Range("S17").Value = 3
$ws.Range("T17").Value = 1.36
# Row 18
$ws.Range("Q18").Value = 2.2
$ws.Range("R18").Value = 1.65
# Row 19
$ws.Range("G19").Value = 2.25
$ws.Range("I19").Value = 3.1
$ws.Range("J19").Value = 3
$ws.Range("L19").Value = 3.75
$ws.Range("AB19").Value = 21
$ws.Range("AC19").Value = 19
$ws.Range("AF19").Value = 6.5
$ws.Range("AI19").Value = 301
$ws.Range("AJ19").Value = 9
$ws.Range("AM19").Value = 34
$ws.Range("AN19").Value = 26
# Row 20
$ws.Range("Q20").Value = 1.65
$ws.Range("R20").Value = 2.2
# Row 21
$ws.Range("G21").Value = 2.5
$ws.Range("H21").Value = 3.3
$ws.Range("I21").Value = 2.75
$ws.Range("J21").Value = 3.2
$ws.Range("Q21").Value = 1.98
$ws.Range("R21").Value = 1.88
$ws.Range("Z21").Value = 12
$ws.Range("AB21").Value = 23
$ws.Range("AF21").Value = 6.5
$ws.Range("AH21").Value = 41
$ws.Range("AM21").Value = 29
# Row 23
$ws.Range("G23").Value = 1.48
$ws.Range("H23").Value = 4.1
$ws.Range("J23").Value = 2.1
$ws.Range("L23").Value = 7
$ws.Range("AA23").Value = 8.5
$ws.Range("AB23").Value = 10
$ws.Range("AE23").Value = 8.5
$ws.Range("AF23").Value = 8
$ws.Range("AJ23").Value = 15
# Row 24
$ws.Range("U24").Value = 1.57
$ws.Range("V24").Value = 2.25
$ws.Range("AM24").Value = 41
$ws.Range("AP24").Value = 1.98
$ws.Range("AQ24").Value = 1.88
# Row 25
$ws.Range("O25").Value = 1.25
$ws.Range("P25").Value = 3.75
$ws.Range("Q25").Value = 1.9
$ws.Range("R25").Value = 1.95
# Row 26
$ws.Range("I26").Value = 2.35
$ws.Range("L26").Value = 3.1
$ws.Range("O26").Value = 1.36
$ws.Range("P26").Value = 3
$ws.Range("Q26").Value = 2.15
$ws.Range("R26").Value = 1.67
$ws.Range("U26").Value = 1.5
$ws.Range("V26").Value = 2.5
$ws.Range("Y26").Value = 8.5
$ws.Range("AC26").Value = 26
$ws.Range("AE26").Value = 8
$ws.Range("AG26").Value = 15
$ws.Range("AI26").Value = 351
$ws.Range("AL26").Value = 9.5
$ws.Range("AM26").Value = 21
# Row 29
$ws.Range("G29").Value = 1.48
$ws.Range("I29").Value = 6.5
$ws.Range("L29").Value = 6
$ws.Range("M29").Value = 1.04
$ws.Range("N29").Value = 13
$ws.Range("Q29").Value = 1.67
$ws.Range("R29").Value = 2.15
$ws.Range("U29").Value = 1.3
$ws.Range("V29").Value = 3.4
$ws.Range("W29").Value = 1.8
$ws.Range("X29").Value = 1.95
# Row 33
$ws.Range("M33").Value = 1.06
$ws.Range("N33").Value = 8
$ws.Range("R33").Value = 1.8
$ws.Range("AC33").Value = 21
$ws.Range("AM33").Value = 23
$ws.Range("AO33").Value = 29
